$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 202, pushing existing rows 202-243 down to 203-244.
$ws.Rows.Item(202).Insert()

# Populate the newly inserted row 202 with the new record.
$ws.Cells.Item(202, 1).Value = 11
$ws.Cells.Item(202, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(202, 3).Value = "Bíobío"
$ws.Cells.Item(202, 4).Value = 44943
$ws.Cells.Item(202, 5).Value = 8
$ws.Cells.Item(202, 6).Value = 100112003
$ws.Cells.Item(202, 7).Value = "Ajo"
$ws.Cells.Item(202, 8).Value = "Chino"
$ws.Cells.Item(202, 9).Value = "Primera"
$ws.Cells.Item(202, 10).Value = 220
$ws.Cells.Item(202, 11).Value = 14000
$ws.Cells.Item(202, 12).Value = 15000
$ws.Cells.Item(202, 13).Value = 14455
$ws.Cells.Item(202, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(202, 15).Value = "China"
$ws.Cells.Item(202, 16).Value = 1446
$ws.Cells.Item(202, 17).Value = 10
$ws.Cells.Item(202, 18).Value = "Hortaliza"
